# Scheduled runner update: refresh market-derived columns H:N (current average
# prices, computed Leve costs/profits) across all job sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 5072.4
$ws.Cells.Item(17, 10).Value = 4199.75
$ws.Cells.Item(17, 12).Value = 12599.25
$ws.Cells.Item(17, 14).Value = -12935.25
$ws.Cells.Item(70, 8).Value = 7885.2856
$ws.Cells.Item(70, 9).Value = 4500
$ws.Cells.Item(70, 10).Value = 9239.4
$ws.Cells.Item(70, 11).Value = 13500
$ws.Cells.Item(70, 12).Value = 27718.2
$ws.Cells.Item(70, 13).Value = -13230
$ws.Cells.Item(70, 14).Value = -28258.2
$ws.Cells.Item(73, 8).Value = 7885.2856
$ws.Cells.Item(73, 9).Value = 4500
$ws.Cells.Item(73, 10).Value = 9239.4
$ws.Cells.Item(73, 11).Value = 13500
$ws.Cells.Item(73, 12).Value = 27718.2
$ws.Cells.Item(73, 13).Value = -12564
$ws.Cells.Item(73, 14).Value = -29590.2
$ws.Cells.Item(111, 8).Value = 450
$ws.Cells.Item(111, 9).Value = 450
$ws.Cells.Item(111, 10).Value = 0
$ws.Cells.Item(111, 11).Value = 1350
$ws.Cells.Item(111, 12).Value = 0
$ws.Cells.Item(111, 13).Value = 1717
$ws.Cells.Item(111, 14).ClearContents()
$ws.Cells.Item(137, 8).Value = 126851
$ws.Cells.Item(137, 9).Value = 0
$ws.Cells.Item(137, 10).Value = 126851
$ws.Cells.Item(137, 11).Value = 0
$ws.Cells.Item(137, 12).Value = 380553
$ws.Cells.Item(137, 13).ClearContents()
$ws.Cells.Item(137, 14).Value = -385653
$ws.Cells.Item(138, 8).Value = 1755.16
$ws.Cells.Item(138, 9).Value = 1187.6552
$ws.Cells.Item(138, 10).Value = 1986.9578
$ws.Cells.Item(138, 11).Value = 3562.9656
$ws.Cells.Item(138, 12).Value = 5960.873399999999
$ws.Cells.Item(138, 13).Value = 1577.0344
$ws.Cells.Item(138, 14).Value = -16240.8734

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1326.21
$ws.Cells.Item(32, 9).Value = 1326.21
$ws.Cells.Item(32, 10).Value = 0
$ws.Cells.Item(32, 11).Value = 1326.21
$ws.Cells.Item(32, 12).Value = 0
$ws.Cells.Item(32, 13).Value = -1039.21
$ws.Cells.Item(32, 14).ClearContents()
$ws.Cells.Item(80, 8).Value = 48666.668
$ws.Cells.Item(80, 10).Value = 48666.668
$ws.Cells.Item(80, 12).Value = 48666.668
$ws.Cells.Item(80, 14).Value = -50662.668
$ws.Cells.Item(83, 8).Value = 48666.668
$ws.Cells.Item(83, 10).Value = 48666.668
$ws.Cells.Item(83, 12).Value = 146000.004
$ws.Cells.Item(83, 14).Value = -155984.004
$ws.Cells.Item(97, 8).Value = 1960
$ws.Cells.Item(97, 9).Value = 2033.3846
$ws.Cells.Item(97, 10).Value = 1801
$ws.Cells.Item(97, 11).Value = 2033.3846
$ws.Cells.Item(97, 12).Value = 1801
$ws.Cells.Item(97, 13).Value = -1537.3846
$ws.Cells.Item(97, 14).Value = -2793
$ws.Cells.Item(102, 8).Value = 1675
$ws.Cells.Item(102, 9).Value = 1675
$ws.Cells.Item(102, 11).Value = 1675
$ws.Cells.Item(102, 13).Value = -53
$ws.Cells.Item(122, 8).Value = 1748.8667
$ws.Cells.Item(122, 9).Value = 1171.7693
$ws.Cells.Item(122, 10).Value = 5500
$ws.Cells.Item(122, 11).Value = 3515.3079
$ws.Cells.Item(122, 12).Value = 16500
$ws.Cells.Item(122, 13).Value = -1065.3079
$ws.Cells.Item(122, 14).Value = -21400
$ws.Cells.Item(132, 8).Value = 2053.0881
$ws.Cells.Item(132, 9).Value = 1666.2273
$ws.Cells.Item(132, 11).Value = 4998.6819
$ws.Cells.Item(132, 13).Value = -2468.6819

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 6915.12
$ws.Cells.Item(134, 9).Value = 9069.412
$ws.Cells.Item(134, 10).Value = 2337.25
$ws.Cells.Item(134, 11).Value = 27208.236
$ws.Cells.Item(134, 12).Value = 7011.75
$ws.Cells.Item(134, 13).Value = -24673.236
$ws.Cells.Item(134, 14).Value = -12081.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 2559517.2
$ws.Cells.Item(58, 9).Value = 3625163.8
$ws.Cells.Item(58, 10).Value = 1965.2
$ws.Cells.Item(58, 11).Value = 3625163.8
$ws.Cells.Item(58, 12).Value = 1965.2
$ws.Cells.Item(58, 13).Value = -3624960.8
$ws.Cells.Item(58, 14).Value = -2371.2
$ws.Cells.Item(59, 8).Value = 19000
$ws.Cells.Item(59, 10).Value = 19000
$ws.Cells.Item(59, 12).Value = 19000
$ws.Cells.Item(59, 14).Value = -21290
$ws.Cells.Item(99, 8).Value = 771839.6
$ws.Cells.Item(99, 9).Value = 1252492.6
$ws.Cells.Item(99, 10).Value = 2794.8
$ws.Cells.Item(99, 11).Value = 1252492.6
$ws.Cells.Item(99, 12).Value = 2794.8
$ws.Cells.Item(99, 13).Value = -1250994.6
$ws.Cells.Item(99, 14).Value = -5790.8
$ws.Cells.Item(126, 8).Value = 771839.6
$ws.Cells.Item(126, 9).Value = 1252492.6
$ws.Cells.Item(126, 10).Value = 2794.8
$ws.Cells.Item(126, 11).Value = 3757477.8
$ws.Cells.Item(126, 12).Value = 8384.400000000001
$ws.Cells.Item(126, 13).Value = -3755007.8
$ws.Cells.Item(126, 14).Value = -13324.4
$ws.Cells.Item(134, 8).Value = 2590.6191
$ws.Cells.Item(134, 10).Value = 2602.7778
$ws.Cells.Item(134, 12).Value = 7808.3334
$ws.Cells.Item(134, 14).Value = -12878.3334
$ws.Cells.Item(136, 8).Value = 2559517.2
$ws.Cells.Item(136, 9).Value = 3625163.8
$ws.Cells.Item(136, 10).Value = 1965.2
$ws.Cells.Item(136, 11).Value = 10875491.4
$ws.Cells.Item(136, 12).Value = 5895.6
$ws.Cells.Item(136, 13).Value = -10872941.4
$ws.Cells.Item(136, 14).Value = -10995.6

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 1328106.8
$ws.Cells.Item(4, 9).Value = 1577056.8
$ws.Cells.Item(4, 11).Value = 4731170.4
$ws.Cells.Item(4, 13).Value = -4731058.4
$ws.Cells.Item(5, 8).Value = 561.2727
$ws.Cells.Item(5, 9).Value = 498.35294
$ws.Cells.Item(5, 11).Value = 1495.05882
$ws.Cells.Item(5, 13).Value = -1383.05882
$ws.Cells.Item(69, 8).Value = 2675.4443
$ws.Cells.Item(69, 9).Value = 2299.8
$ws.Cells.Item(69, 11).Value = 6899.400000000001
$ws.Cells.Item(69, 13).Value = -6088.400000000001
$ws.Cells.Item(72, 8).Value = 2675.4443
$ws.Cells.Item(72, 9).Value = 2299.8
$ws.Cells.Item(72, 11).Value = 20698.2
$ws.Cells.Item(72, 13).Value = -16642.2
$ws.Cells.Item(97, 8).Value = 910.25
$ws.Cells.Item(97, 9).Value = 280.33334
$ws.Cells.Item(97, 11).Value = 841.0000200000001
$ws.Cells.Item(97, 13).Value = -345.0000200000001
$ws.Cells.Item(113, 8).Value = 43188.92
$ws.Cells.Item(113, 10).Value = 988.4211
$ws.Cells.Item(113, 12).Value = 2965.2633
$ws.Cells.Item(113, 14).Value = -7305.263300000001
$ws.Cells.Item(135, 8).Value = 561.2727
$ws.Cells.Item(135, 9).Value = 498.35294
$ws.Cells.Item(135, 11).Value = 4485.17646
$ws.Cells.Item(135, 13).Value = -1950.17646
$ws.Cells.Item(137, 8).Value = 3913.75
$ws.Cells.Item(137, 9).Value = 1231.5
$ws.Cells.Item(137, 11).Value = 3694.5
$ws.Cells.Item(137, 13).Value = 1405.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 2602.75
$ws.Cells.Item(80, 9).Value = 2591.5293
$ws.Cells.Item(80, 10).Value = 2666.3333
$ws.Cells.Item(80, 11).Value = 2591.5293
$ws.Cells.Item(80, 12).Value = 2666.3333
$ws.Cells.Item(80, 13).Value = -1593.5293
$ws.Cells.Item(80, 14).Value = -4662.3333
$ws.Cells.Item(83, 8).Value = 2602.75
$ws.Cells.Item(83, 9).Value = 2591.5293
$ws.Cells.Item(83, 10).Value = 2666.3333
$ws.Cells.Item(83, 11).Value = 12957.6465
$ws.Cells.Item(83, 12).Value = 13331.6665
$ws.Cells.Item(83, 13).Value = -7965.646500000001
$ws.Cells.Item(83, 14).Value = -23315.6665
$ws.Cells.Item(97, 8).Value = 978.34375
$ws.Cells.Item(97, 9).Value = 972.84
$ws.Cells.Item(97, 11).Value = 972.84
$ws.Cells.Item(97, 13).Value = -476.84
$ws.Cells.Item(132, 8).Value = 1482299.4
$ws.Cells.Item(132, 9).Value = 1925663
$ws.Cells.Item(132, 10).Value = 4420.8335
$ws.Cells.Item(132, 11).Value = 5776989
$ws.Cells.Item(132, 12).Value = 13262.5005
$ws.Cells.Item(132, 13).Value = -5774459
$ws.Cells.Item(132, 14).Value = -18322.5005

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 5756.3335
$ws.Cells.Item(16, 9).Value = 6238.727
$ws.Cells.Item(16, 11).Value = 6238.727
$ws.Cells.Item(16, 13).Value = -6068.727
$ws.Cells.Item(22, 8).Value = 1161.96
$ws.Cells.Item(22, 9).Value = 654.2
$ws.Cells.Item(22, 10).Value = 1500.4667
$ws.Cells.Item(22, 11).Value = 654.2
$ws.Cells.Item(22, 12).Value = 1500.4667
$ws.Cells.Item(22, 13).Value = -359.2
$ws.Cells.Item(22, 14).Value = -2090.4667
$ws.Cells.Item(27, 8).Value = 1161.96
$ws.Cells.Item(27, 9).Value = 654.2
$ws.Cells.Item(27, 10).Value = 1500.4667
$ws.Cells.Item(27, 11).Value = 654.2
$ws.Cells.Item(27, 12).Value = 1500.4667
$ws.Cells.Item(27, 13).Value = -547.2
$ws.Cells.Item(27, 14).Value = -1714.4667
$ws.Cells.Item(46, 8).Value = 1980.8889
$ws.Cells.Item(46, 9).Value = 1648.3334
$ws.Cells.Item(46, 10).Value = 2147.1667
$ws.Cells.Item(46, 11).Value = 1648.3334
$ws.Cells.Item(46, 12).Value = 2147.1667
$ws.Cells.Item(46, 13).Value = -1460.3334
$ws.Cells.Item(46, 14).Value = -2523.1667
$ws.Cells.Item(122, 8).Value = 2850
$ws.Cells.Item(122, 9).Value = 2757.1428
$ws.Cells.Item(122, 11).Value = 8271.428400000001
$ws.Cells.Item(122, 13).Value = -5821.428400000001
$ws.Cells.Item(132, 8).Value = 2240.6758
$ws.Cells.Item(132, 9).Value = 1829.0714
$ws.Cells.Item(132, 11).Value = 5487.2142
$ws.Cells.Item(132, 13).Value = -2957.2142

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 646.7692
$ws.Cells.Item(113, 9).Value = 494.6
$ws.Cells.Item(113, 10).Value = 741.875
$ws.Cells.Item(113, 11).Value = 1483.8
$ws.Cells.Item(113, 12).Value = 2225.625
$ws.Cells.Item(113, 13).Value = 686.1999999999998
$ws.Cells.Item(113, 14).Value = -6565.625
$ws.Cells.Item(122, 8).Value = 27123.033
$ws.Cells.Item(122, 9).Value = 34454.332
$ws.Cells.Item(122, 11).Value = 103362.996
$ws.Cells.Item(122, 13).Value = -100912.996
$ws.Cells.Item(132, 8).Value = 1259.6383
$ws.Cells.Item(132, 9).Value = 1095.95
$ws.Cells.Item(132, 10).Value = 2195
$ws.Cells.Item(132, 11).Value = 3287.85
$ws.Cells.Item(132, 12).Value = 6585
$ws.Cells.Item(132, 13).Value = -757.8500000000004
$ws.Cells.Item(132, 14).Value = -11645
$ws.Cells.Item(136, 8).Value = 26456278
$ws.Cells.Item(136, 9).Value = 42735570
$ws.Cells.Item(136, 10).Value = 2433.5
$ws.Cells.Item(136, 11).Value = 128206710
$ws.Cells.Item(136, 12).Value = 7300.5
$ws.Cells.Item(136, 13).Value = -128204160
$ws.Cells.Item(136, 14).Value = -12400.5
